{"js": "const replacements = [\n  [\n    \"2. System exibe a listagem dos Avaliacoes cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' \",\n    \"2. System exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' \"\n  ],\n  [\n    \"12. System apresenta em 'Metas' os campos 'Nivel' preenchido corretamente \",\n    \"12. System apresenta em 'Metas' os campos 'Nivel' preenchidos corretamente \"\n  ],\n  [\n    \"4. System apresenta o formulario para e alteracao de Avaliacao \",\n    \"4. System apresenta o formulario para e alteracao da Avaliacao \"\n  ],\n  [\n    \"5. Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' somente leitura bs 11\",\n    \"5. Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' estao em modo somente leitura bs 11\"\n  ],\n  [\n    \"AF[2] \\u2013 Confirmar Exclusao do Avaliacao\",\n    \"AF[2] \\u2013 Confirmar Exclusao da Avaliacao\"\n  ],\n  [\n    \"6. System exibe a listagem dos Avaliacoes sem o Avaliacao excluido ef[3,4]\",\n    \"6. System exibe a listagem das Avaliacoes sem a Avaliacao excluida ef[3,4]\"\n  ],\n  [\n    \"AF[3] \\u2013 Negar Exclusao do Avaliacao\",\n    \"AF[3] \\u2013 Nao Confirmar a Exclusao da Avaliacao\"\n  ],\n  [\n    \"6. System exibe a listagem dos Avaliacoes com o Avaliacao excluido \",\n    \"6. System exibe a listagem das Avaliacoes com a Avaliacao nao excluida \"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // Search for the exact (whole) run text and replace just that range in\n  // place, so surrounding runs (bookmarks, empty runs, run formatting)\n  // are left completely untouched.\n  const searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"2. System exibe a listagem dos Avaliacoes cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' \"\n        New = \"2. System exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' \"\n    },\n    @{\n        Old = \"12. System apresenta em 'Metas' os campos 'Nivel' preenchido corretamente \"\n        New = \"12. System apresenta em 'Metas' os campos 'Nivel' preenchidos corretamente \"\n    },\n    @{\n        Old = \"4. System apresenta o formulario para e alteracao de Avaliacao \"\n        New = \"4. System apresenta o formulario para e alteracao da Avaliacao \"\n    },\n    @{\n        Old = \"5. Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' somente leitura bs 11\"\n        New = \"5. Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' estao em modo somente leitura bs 11\"\n    },\n    @{\n        Old = \"AF[2] \u2013 Confirmar Exclusao do Avaliacao\"\n        New = \"AF[2] \u2013 Confirmar Exclusao da Avaliacao\"\n    },\n    @{\n        Old = \"6. System exibe a listagem dos Avaliacoes sem o Avaliacao excluido ef[3,4]\"\n        New = \"6. System exibe a listagem das Avaliacoes sem a Avaliacao excluida ef[3,4]\"\n    },\n    @{\n        Old = \"AF[3] \u2013 Negar Exclusao do Avaliacao\"\n        New = \"AF[3] \u2013 Nao Confirmar a Exclusao da Avaliacao\"\n    },\n    @{\n        Old = \"6. System exibe a listagem dos Avaliacoes com o Avaliacao excluido \"\n        New = \"6. System exibe a listagem das Avaliacoes com a Avaliacao nao excluida \"\n    }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if ($found) {\n        # Assign directly to Range.Text (instead of Find.Replacement) so\n        # Word's AutoCorrect / smart-quotes AutoFormat does not mangle the\n        # straight apostrophes in the replacement text.\n        $range.Text = $r.New\n    }\n}\n"}
